# Auto-generated Excel COM-interop script applying the Shiva_Profits.xlsx diff
# (FFXIV leve-profit calcs refreshed by the scheduled market-data runner).
$wb = $excel.ActiveWorkbook

# --- ALC!row32 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 50012064
$ws.Cells.Item(32, 9).Value = 66681668
$ws.Cells.Item(32, 11).Value = 66681668
$ws.Cells.Item(32, 13).Value = -66681342

# --- ALC!row57 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(57, 8).Value = 142412.25
$ws.Cells.Item(57, 10).Value = 142412.25
$ws.Cells.Item(57, 12).Value = 427236.75
$ws.Cells.Item(57, 14).Value = -428234.75

# --- ALC!row98 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1348.7222
$ws.Cells.Item(98, 9).Value = 1284.9656
$ws.Cells.Item(98, 10).Value = 1612.8572
$ws.Cells.Item(98, 11).Value = 1284.9656
$ws.Cells.Item(98, 12).Value = 1612.8572
$ws.Cells.Item(98, 13).Value = 213.0344
$ws.Cells.Item(98, 14).Value = -4608.8572

# --- ALC!row106 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 30308888
$ws.Cells.Item(106, 9).Value = 33338778
$ws.Cells.Item(106, 11).Value = 33338778
$ws.Cells.Item(106, 13).Value = -33338147

# --- ALC!row122 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 1348.7222
$ws.Cells.Item(122, 9).Value = 1284.9656
$ws.Cells.Item(122, 10).Value = 1612.8572
$ws.Cells.Item(122, 11).Value = 3854.8968
$ws.Cells.Item(122, 12).Value = 4838.571599999999
$ws.Cells.Item(122, 13).Value = -1404.8968
$ws.Cells.Item(122, 14).Value = -9738.571599999999

# --- ARM!row32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9398.134
$ws.Cells.Item(32, 9).Value = 8832.5
$ws.Cells.Item(32, 10).Value = 13074.75
$ws.Cells.Item(32, 11).Value = 8832.5
$ws.Cells.Item(32, 12).Value = 13074.75
$ws.Cells.Item(32, 13).Value = -8545.5
$ws.Cells.Item(32, 14).Value = -13648.75

# --- ARM!row37 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(37, 8).Value = 31011.334

# --- ARM!row74 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 6865.467
$ws.Cells.Item(74, 9).Value = 6865.467
$ws.Cells.Item(74, 11).Value = 6865.467
$ws.Cells.Item(74, 13).Value = -5991.467

# --- ARM!row77 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 6865.467
$ws.Cells.Item(77, 9).Value = 6865.467
$ws.Cells.Item(77, 11).Value = 34327.335
$ws.Cells.Item(77, 13).Value = -29959.335

# --- ARM!row94 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(94, 8).Value = 255000
$ws.Cells.Item(94, 9).Value = 10000
$ws.Cells.Item(94, 10).Value = 500000
$ws.Cells.Item(94, 11).Value = 10000
$ws.Cells.Item(94, 12).Value = 500000
$ws.Cells.Item(94, 13).Value = -9099
$ws.Cells.Item(94, 14).Value = -501802

# --- ARM!row109 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 14).Value = 0
$ws.Cells.Item(109, 12).ClearContents()

# --- ARM!row119 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(119, 8).Value = 41997.332
$ws.Cells.Item(119, 10).Value = 41997.332
$ws.Cells.Item(119, 12).Value = 41997.332
$ws.Cells.Item(119, 14).Value = -51673.332

# --- ARM!row132 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3830.1035
$ws.Cells.Item(132, 9).Value = 3487.926
$ws.Cells.Item(132, 10).Value = 8449.5
$ws.Cells.Item(132, 11).Value = 10463.778
$ws.Cells.Item(132, 12).Value = 25348.5
$ws.Cells.Item(132, 13).Value = -7933.778
$ws.Cells.Item(132, 14).Value = -30408.5

# --- ARM!row141 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(141, 8).Value = 83344
$ws.Cells.Item(141, 10).Value = 83344
$ws.Cells.Item(141, 12).Value = 83344
$ws.Cells.Item(141, 14).Value = -93704

# --- BSM!row94 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2750.5454
$ws.Cells.Item(94, 9).Value = 2760.3333
$ws.Cells.Item(94, 10).Value = 2738.8
$ws.Cells.Item(94, 11).Value = 2760.3333
$ws.Cells.Item(94, 12).Value = 2738.8
$ws.Cells.Item(94, 13).Value = -2309.3333
$ws.Cells.Item(94, 14).Value = -3640.8

# --- BSM!row107 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1439.7
$ws.Cells.Item(107, 9).Value = 1249.92
$ws.Cells.Item(107, 11).Value = 1249.92
$ws.Cells.Item(107, 13).Value = 670.0799999999999

# --- BSM!row134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2706.7288
$ws.Cells.Item(134, 9).Value = 2647.9363
$ws.Cells.Item(134, 11).Value = 7943.8089
$ws.Cells.Item(134, 13).Value = -5408.8089

# --- CRP!row16 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2403.4614
$ws.Cells.Item(16, 9).Value = 1792.6666
$ws.Cells.Item(16, 10).Value = 2927
$ws.Cells.Item(16, 11).Value = 1792.6666
$ws.Cells.Item(16, 12).Value = 2927
$ws.Cells.Item(16, 13).Value = -1505.6666
$ws.Cells.Item(16, 14).Value = -3501

# --- CRP!row31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2398.25
$ws.Cells.Item(31, 9).Value = 2129.3845
$ws.Cells.Item(31, 10).Value = 3563.3333
$ws.Cells.Item(31, 11).Value = 2129.3845
$ws.Cells.Item(31, 12).Value = 3563.3333
$ws.Cells.Item(31, 13).Value = -1834.3845
$ws.Cells.Item(31, 14).Value = -4153.3333

# --- CRP!row34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2398.25
$ws.Cells.Item(34, 9).Value = 2129.3845
$ws.Cells.Item(34, 10).Value = 3563.3333
$ws.Cells.Item(34, 11).Value = 2129.3845
$ws.Cells.Item(34, 12).Value = 3563.3333
$ws.Cells.Item(34, 13).Value = -1927.3845
$ws.Cells.Item(34, 14).Value = -3967.3333

# --- CRP!row113 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 2403.4614
$ws.Cells.Item(113, 9).Value = 1792.6666
$ws.Cells.Item(113, 10).Value = 2927
$ws.Cells.Item(113, 11).Value = 1792.6666
$ws.Cells.Item(113, 12).Value = 2927
$ws.Cells.Item(113, 13).Value = 377.3334
$ws.Cells.Item(113, 14).Value = -7267

# --- CRP!row134 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 3854.6667
$ws.Cells.Item(134, 9).Value = 3373.077
$ws.Cells.Item(134, 11).Value = 10119.231
$ws.Cells.Item(134, 13).Value = -7584.231

# --- CUL!row4 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 75170.71000000001
$ws.Cells.Item(4, 10).Value = 4365.8335
$ws.Cells.Item(4, 12).Value = 13097.5005
$ws.Cells.Item(4, 14).Value = -13321.5005

# --- CUL!row7 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 97.375
$ws.Cells.Item(7, 9).Value = 40
$ws.Cells.Item(7, 10).Value = 499
$ws.Cells.Item(7, 11).Value = 120
$ws.Cells.Item(7, 12).Value = 1497
$ws.Cells.Item(7, 13).Value = -8
$ws.Cells.Item(7, 14).Value = -1721

# --- CUL!row69 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 947
$ws.Cells.Item(69, 9).Value = 947
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 2841
$ws.Cells.Item(69, 14).Value = 0
$ws.Cells.Item(69, 13).Value = -2030
$ws.Cells.Item(69, 12).ClearContents()

# --- CUL!row72 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 8).Value = 947
$ws.Cells.Item(72, 9).Value = 947
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 8523
$ws.Cells.Item(72, 14).Value = 0
$ws.Cells.Item(72, 13).Value = -4467
$ws.Cells.Item(72, 12).ClearContents()

# --- CUL!row80 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 5500
$ws.Cells.Item(80, 10).Value = 5500
$ws.Cells.Item(80, 12).Value = 16500
$ws.Cells.Item(80, 14).Value = -18372

# --- CUL!row81 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 1122
$ws.Cells.Item(81, 9).Value = 1122
$ws.Cells.Item(81, 11).Value = 3366
$ws.Cells.Item(81, 13).Value = -2243

# --- CUL!row83 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value = 5500
$ws.Cells.Item(83, 10).Value = 5500
$ws.Cells.Item(83, 12).Value = 49500
$ws.Cells.Item(83, 14).Value = -58860

# --- CUL!row84 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(84, 8).Value = 1122
$ws.Cells.Item(84, 9).Value = 1122
$ws.Cells.Item(84, 11).Value = 10098
$ws.Cells.Item(84, 13).Value = -4482

# --- CUL!row87 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 14).Value = 0
$ws.Cells.Item(87, 12).ClearContents()

# --- CUL!row90 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 14).Value = 0
$ws.Cells.Item(90, 12).ClearContents()

# --- CUL!row113 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1055.0869
$ws.Cells.Item(113, 9).Value = 942.375
$ws.Cells.Item(113, 10).Value = 1115.2
$ws.Cells.Item(113, 11).Value = 2827.125
$ws.Cells.Item(113, 12).Value = 3345.6
$ws.Cells.Item(113, 13).Value = -657.125
$ws.Cells.Item(113, 14).Value = -7685.6

# --- CUL!row131 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1338300.9
$ws.Cells.Item(131, 9).Value = 2942265.8
$ws.Cells.Item(131, 10).Value = 1663.3334
$ws.Cells.Item(131, 11).Value = 8826797.399999999
$ws.Cells.Item(131, 12).Value = 4990.0002
$ws.Cells.Item(131, 13).Value = -8821757.399999999
$ws.Cells.Item(131, 14).Value = -15070.0002

# --- CUL!row132 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 4239.636
$ws.Cells.Item(132, 9).Value = 4083.5557
$ws.Cells.Item(132, 10).Value = 4942
$ws.Cells.Item(132, 11).Value = 36752.0013
$ws.Cells.Item(132, 12).Value = 44478
$ws.Cells.Item(132, 13).Value = -34222.0013
$ws.Cells.Item(132, 14).Value = -49538

# --- GSM!row39 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 19000
$ws.Cells.Item(39, 10).Value = 19000
$ws.Cells.Item(39, 12).Value = 19000
$ws.Cells.Item(39, 14).Value = -20064

# --- GSM!row97 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 898.51514
$ws.Cells.Item(97, 9).Value = 829.2222
$ws.Cells.Item(97, 11).Value = 829.2222
$ws.Cells.Item(97, 13).Value = -333.2222

# --- GSM!row122 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3948.2856
$ws.Cells.Item(122, 9).Value = 3536.3103
$ws.Cells.Item(122, 11).Value = 10608.9309
$ws.Cells.Item(122, 13).Value = -8158.930899999999

# --- LTW!row22 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 913.25
$ws.Cells.Item(22, 9).Value = 939.8333
$ws.Cells.Item(22, 10).Value = 833.5
$ws.Cells.Item(22, 11).Value = 939.8333
$ws.Cells.Item(22, 12).Value = 833.5
$ws.Cells.Item(22, 13).Value = -644.8333
$ws.Cells.Item(22, 14).Value = -1423.5

# --- LTW!row27 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 913.25
$ws.Cells.Item(27, 9).Value = 939.8333
$ws.Cells.Item(27, 10).Value = 833.5
$ws.Cells.Item(27, 11).Value = 939.8333
$ws.Cells.Item(27, 12).Value = 833.5
$ws.Cells.Item(27, 13).Value = -832.8333
$ws.Cells.Item(27, 14).Value = -1047.5

# --- LTW!row136 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 11448.971
$ws.Cells.Item(136, 9).Value = 8065.2334
$ws.Cells.Item(136, 11).Value = 24195.7002
$ws.Cells.Item(136, 13).Value = -21645.7002

# --- WVR!row103 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(103, 8).Value = 33400.25
$ws.Cells.Item(103, 10).Value = 33400.25
$ws.Cells.Item(103, 12).Value = 33400.25
$ws.Cells.Item(103, 14).Value = -35744.25
